$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "CIRUGIA PROPUESTA: CIRUGÌA PARA EXERESIS DE COLESTEATOMA"
#   -> "CIRUGIA PROPUESTA: CIRUGÍA PARA EXÉRESIS DE COLESTEATOMA"
# (fix the accent typos: Ì -> Í, EXERESIS -> EXÉRESIS). Formatting (bold,
# Book Antiqua, 20) is unchanged for the whole run.
# ---------------------------------------------------------------------------
$find1 = $d.Content.Duplicate
$found1 = $find1.Find.Execute("CIRUGÌA PARA EXERESIS DE COLESTEATOMA", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "CIRUGÍA PARA EXÉRESIS DE COLESTEATOMA", 2)

# ---------------------------------------------------------------------------
# Change 2: split the "Este documento informativo..." paragraph so that the
# phrase "CIRUGÍA DEL COLESTEATOMA" becomes bold, while the rest of the
# sentence keeps its original (non-bold) formatting.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -like "Este documento informativo pretende explicar*") {
        $sub = $p.Range.Duplicate
        $sub.Find.Execute("CIRUGÍA DEL COLESTEATOMA")
        if ($sub.Find.Found) {
            $sub.Font.Bold = 1
            $sub.Font.Name = "Book Antiqua"
        }
    }
}
